# adds AustralianJPS to notes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (78) with the Australian Journal of Political Science entry.
$ws.Cells.Item(78, 1).Value = "Australian Journal of Political Science"
$ws.Cells.Item(78, 2).Value = "<a href='https://www.tandfonline.com/action/authorSubmission?show=instructions&journalCode=cajp20'target='_blank'>Research Note</a>"
$ws.Cells.Item(78, 3).Value = "4k words"
$ws.Cells.Item(78, 4).Value = 17

# Move/restore the on-screen selection to the new entry, as the author left it.
[void]$ws.Range("G78").Select()
